$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the empty B7/C7 inline-string placeholders
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = ""

# Add the new row 8 data
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Vinicius"
$ws.Range("C8").Value = "'09/04/2022"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "01:00"
$ws.Range("E8").Value = 1000
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 300
$ws.Range("H8").Value = 14000
$ws.Range("J8").Value = "90000 kg"
$ws.Range("K8").Value = "Pendente"
